$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "28.362.13"
Set-TextCell $ws.Range("E2") "  +2.85%  "
Set-TextCell $ws.Range("D3") "1.870.07"
Set-TextCell $ws.Range("E3") "  +1.23%  "
Set-TextCell $ws.Range("D5") "338.99"
Set-TextCell $ws.Range("E5") "  +2.02%  "
Set-TextCell $ws.Range("E6") "  -0.18%  "
Set-TextCell $ws.Range("D7") "0.4703"
Set-TextCell $ws.Range("E7") "  +1.66%  "
Set-TextCell $ws.Range("D8") "0.3947"
Set-TextCell $ws.Range("E8") "  +2.50%  "
Set-TextCell $ws.Range("D9") "47.32"
Set-TextCell $ws.Range("E9") "  +2.91%  "
Set-TextCell $ws.Range("D10") "0.08000"
Set-TextCell $ws.Range("E10") "  +1.15%  "
Set-TextCell $ws.Range("D11") "1.008"
Set-TextCell $ws.Range("E11") "  +1.52%  "
Set-TextCell $ws.Range("D12") "21.95"
Set-TextCell $ws.Range("E12") "  +2.35%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D13") "6.008"
Set-TextCell $ws.Range("E13") "  +1.72%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D14") "1.868.98"
Set-TextCell $ws.Range("E14") "  +0.64%  "
Set-TextCell $ws.Range("D15") "7.284"
Set-TextCell $ws.Range("E15") "  +2.75%  "
Set-TextCell $ws.Range("D16") "91.24"
Set-TextCell $ws.Range("E16") "  +2.91%  "
Set-TextCell $ws.Range("D17") "1.003"
Set-TextCell $ws.Range("E17") "  -0.13%  "
Set-TextCell $ws.Range("D18") "0.00001042"
Set-TextCell $ws.Range("E18") "  +0.80%  "
Set-TextCell $ws.Range("D19") "0.06595"
Set-TextCell $ws.Range("E19") "  -0.88%  "
Set-TextCell $ws.Range("D20") "17.71"
Set-TextCell $ws.Range("E20") "  +3.93%  "
Set-TextCell $ws.Range("D21") "1.000"
Set-TextCell $ws.Range("E21") "  -0.14%  "
Set-TextCell $ws.Range("D22") "28.349.79"
Set-TextCell $ws.Range("E22") "  +2.79%  "
Set-TextCell $ws.Range("E23") "  +1.65%  "
Set-TextCell $ws.Range("E24") "  +1.51%  "
Set-TextCell $ws.Range("D25") "2.287"
Set-TextCell $ws.Range("E25") "  -0.78%  "
Set-TextCell $ws.Range("D26") "2.106.77"
Set-TextCell $ws.Range("E26") "  +1.68%  "
Set-TextCell $ws.Range("D27") "159.79"
Set-TextCell $ws.Range("E27") "  +1.25%  "
Set-TextCell $ws.Range("D28") "19.87"
Set-TextCell $ws.Range("E28") "  +2.05%  "
Set-TextCell $ws.Range("D29") "2.148"
Set-TextCell $ws.Range("E29") "  +2.89%  "
Set-TextCell $ws.Range("D30") "5.497"
Set-TextCell $ws.Range("E30") "  +1.94%  "
Set-TextCell $ws.Range("D31") "120.40"
Set-TextCell $ws.Range("E31") "  +0.64%  "
Set-TextCell $ws.Range("D32") "0.9757"
Set-TextCell $ws.Range("D33") "0.09516"
Set-TextCell $ws.Range("E33") "  +1.41%  "
Set-TextCell $ws.Range("D34") "3.592"
Set-TextCell $ws.Range("E34") "  +0.39%  "
Set-TextCell $ws.Range("D35") "1.382"
Set-TextCell $ws.Range("E35") "  +2.87%  "
Set-TextCell $ws.Range("D36") "5.359"
Set-TextCell $ws.Range("E36") "  +1.63%  "
Set-TextCell $ws.Range("D37") "0.02276"
Set-TextCell $ws.Range("E37") "  +2.64%  "
Set-TextCell $ws.Range("D38") "0.06101"
Set-TextCell $ws.Range("E38") "  +1.68%  "
Set-TextCell $ws.Range("D39") "8.483"
Set-TextCell $ws.Range("E39") "  +2.55%  "
Set-TextCell $ws.Range("E40") "  +0.21%  "
Set-TextCell $ws.Range("D41") "0.5980"
Set-TextCell $ws.Range("E41") "  +1.71%  "
Set-TextCell $ws.Range("E42") "  -0.13%  "
Set-TextCell $ws.Range("D43") "0.1884"
Set-TextCell $ws.Range("E43") "  +1.33%  "
Set-TextCell $ws.Range("D44") "10.39"
Set-TextCell $ws.Range("E44") "  +1.22%  "
Set-TextCell $ws.Range("D45") "1.289"
Set-TextCell $ws.Range("E45") "  +3.69%  "
Set-TextCell $ws.Range("D46") "0.5629"
Set-TextCell $ws.Range("E46") "  +1.00%  "
Set-TextCell $ws.Range("D47") "12.14"
Set-TextCell $ws.Range("E47") "  -0.21%  "
Set-TextCell $ws.Range("D48") "1.966"
Set-TextCell $ws.Range("E48") "  +3.94%  "
Set-TextCell $ws.Range("D49") "0.06925"
Set-TextCell $ws.Range("E49") "  +3.66%  "
Set-TextCell $ws.Range("D50") "111.37"
Set-TextCell $ws.Range("E50") "  +0.79%  "
Set-TextCell $ws.Range("D51") "2.022"
Set-TextCell $ws.Range("E51") "  +13.47%  "
